# Updated cryptos list values (Price + Volume(1h)) per the Nov 9 2024 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '75.970.02'
$c.ClearFormats()
$ws.Range("E2").Value = '  -0.02%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.014.32'
$c.ClearFormats()
$ws.Range("E3").Value = '  +4.09%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("E5").Value = '  -0.25%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '615.42'
$c.ClearFormats()
$ws.Range("E6").Value = '  +4.43%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  +0.92%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.205'
$c.ClearFormats()
$ws.Range("E9").Value = '  +7.02%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '3.011.45'
$c.ClearFormats()
$ws.Range("E10").Value = '  +3.99%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.436'
$c.ClearFormats()
$ws.Range("E11").Value = '  -0.58%  '

$ws.Range("E12").Value = '  -0.19%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.20'
$c.ClearFormats()
$ws.Range("E13").Value = '  +7.46%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '3.570.02'
$c.ClearFormats()
$ws.Range("E14").Value = '  +4.15%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '28.79'
$c.ClearFormats()
$ws.Range("E15").Value = '  +4.33%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '75.872.85'
$c.ClearFormats()
$ws.Range("E16").Value = '  -0.06%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.0000191'
$c.ClearFormats()
$ws.Range("E17").Value = '  +3.20%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.004.78'
$c.ClearFormats()
$ws.Range("E18").Value = '  +3.06%  '

$ws.Range("E19").Value = '  +3.07%  '

$ws.Range("E20").Value = '  +3.90%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '377.24'
$c.ClearFormats()
$ws.Range("E21").Value = '  +3.28%  '

$ws.Range("E22").Value = '  +6.25%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.36'
$c.ClearFormats()
$ws.Range("E23").Value = '  +2.58%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '3.171.27'
$c.ClearFormats()
$ws.Range("E24").Value = '  +4.12%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '72.15'
$c.ClearFormats()
$ws.Range("E25").Value = '  +0.97%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("E27").Value = '  +3.39%  '

$ws.Range("E28").Value = '  +2.87%  '

$ws.Range("E29").Value = '  +3.75%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.ClearFormats()
$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("E31").Value = '  +3.61%  '

$ws.Range("E32").Value = '  +2.56%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '491.14'
$c.ClearFormats()
$ws.Range("E33").Value = '  +0.28%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.90'
$c.ClearFormats()
$ws.Range("E34").Value = '  +5.75%  '

$ws.Range("E35").Value = '  +0.06%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '20.49'
$c.ClearFormats()
$ws.Range("E36").Value = '  +3.17%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.121'
$c.ClearFormats()
$ws.Range("E37").Value = '  +12.08%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '161.97'
$c.ClearFormats()
$ws.Range("E38").Value = '  -2.24%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '20.03'
$c.ClearFormats()
$ws.Range("E39").Value = '  +1.74%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '190.18'
$c.ClearFormats()
$ws.Range("E40").Value = '  +7.36%  '

$ws.Range("E41").Value = '  -1.74%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.103'
$c.ClearFormats()
$ws.Range("E42").Value = '  -3.97%  '

$ws.Range("E43").Value = '  -0.03%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.09'
$c.ClearFormats()
$ws.Range("E44").Value = '  +5.87%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.769'
$c.ClearFormats()
$ws.Range("E45").Value = '  +19.07%  '

$ws.Range("E46").Value = '  +7.38%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '41.15'
$c.ClearFormats()
$ws.Range("E47").Value = '  +2.60%  '

$ws.Range("E48").Value = '  +1.13%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.ClearFormats()
$ws.Range("E49").Value = '  +9.21%  '

$ws.Range("E50").Value = '  +3.40%  '

$ws.Range("E51").Value = '  +1.53%  '
